# Edit: Slide 7 ("NavMesh Example") - renewal of the Bake instructions text,
# resize of the containing textbox, and removal of stale local-path picture
# descriptions (descr attrs) left over from PolarisOffice round-tripping.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# --- 1. Rebuild the instructional paragraph text ("12. ...") -------------
# Shape 2 ("Rect 0") holds: "12." (bold/blue) + " " (plain) + the long
# Korean sentence split across many runs. We keep the "12." and " " runs
# untouched and replace the remainder with a single merged run/sentence,
# matching how the text was rewritten by hand.
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$len = $tr.Length
$rest = $tr.Characters(5, $len - 4)
[void]$rest.Delete()
[void]$tr.InsertAfter("그런 다음 NavMeshAgent가 이동할 수 있는 범위를 설정하였으면 Bake를 선택합니다.")

# Restore the shape's box size (text edits above can perturb the auto-fit
# height) and apply the slightly wider target width from the diff.
$sh.Width = 337.300004
$sh.Height = 72.75

# --- 2. Drop the leftover local-machine descriptions on the two pictures -
$pic1 = $s.Shapes.Item(3)
$pic1.AlternativeText = ""

$pic2 = $s.Shapes.Item(5)
$pic2.AlternativeText = ""
